$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, $value)
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.ClearFormats()
}

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D14") "7.147"
$ws.Range("E14").Value = "  -5.97%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D15") "5.859"
$ws.Range("E15").Value = "  -5.36%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D48") "1.917"
$ws.Range("E48").Value = "  -6.94%  "

$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D49") "3.349"
$ws.Range("E49").Value = "  -1.81%  "

Set-TextValue $ws.Range("D2") "27.259.89"
$ws.Range("E2").Value = "  -4.45%  "

Set-TextValue $ws.Range("D3") "1.853.78"
$ws.Range("E3").Value = "  -5.69%  "

$ws.Range("E4").Value = "  -1.22%  "

Set-TextValue $ws.Range("D5") "320.76"
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("E6").Value = "  -1.11%  "

Set-TextValue $ws.Range("D7") "0.4493"
$ws.Range("E7").Value = "  -5.69%  "

Set-TextValue $ws.Range("D8") "0.3843"
$ws.Range("E8").Value = "  -5.04%  "

Set-TextValue $ws.Range("D9") "47.71"
$ws.Range("E9").Value = "  -11.76%  "

Set-TextValue $ws.Range("D10") "0.07867"
$ws.Range("E10").Value = "  -7.16%  "

Set-TextValue $ws.Range("D11") "1.015"
$ws.Range("E11").Value = "  -4.28%  "

Set-TextValue $ws.Range("D12") "21.31"
$ws.Range("E12").Value = "  -4.87%  "

Set-TextValue $ws.Range("D13") "1.858.53"
$ws.Range("E13").Value = "  -8.28%  "

Set-TextValue $ws.Range("D16") "1.001"
$ws.Range("E16").Value = "  -1.28%  "

Set-TextValue $ws.Range("D17") "0.00001028"
$ws.Range("E17").Value = "  -3.80%  "

Set-TextValue $ws.Range("D18") "85.53"
$ws.Range("E18").Value = "  -5.87%  "

Set-TextValue $ws.Range("D19") "0.06527"
$ws.Range("E19").Value = "  -1.75%  "

Set-TextValue $ws.Range("D20") "16.94"
$ws.Range("E20").Value = "  -8.57%  "

Set-TextValue $ws.Range("D21") "1.000"
$ws.Range("E21").Value = "  -1.21%  "

Set-TextValue $ws.Range("D22") "5.483"
$ws.Range("E22").Value = "  -6.46%  "

Set-TextValue $ws.Range("D23") "27.267.18"
$ws.Range("E23").Value = "  -4.56%  "

Set-TextValue $ws.Range("D24") "10.76"
$ws.Range("E24").Value = "  -6.19%  "

Set-TextValue $ws.Range("D25") "2.262"
$ws.Range("E25").Value = "  -1.72%  "

Set-TextValue $ws.Range("D26") "2.084.93"
$ws.Range("E26").Value = "  -7.50%  "

Set-TextValue $ws.Range("D27") "151.51"
$ws.Range("E27").Value = "  -2.76%  "

Set-TextValue $ws.Range("D28") "19.66"
$ws.Range("E28").Value = "  -3.39%  "

Set-TextValue $ws.Range("D29") "2.061"
$ws.Range("E29").Value = "  -4.91%  "

Set-TextValue $ws.Range("D30") "5.457"
$ws.Range("E30").Value = "  -7.45%  "

$ws.Range("E31").Value = "  -3.64%  "

$ws.Range("E32").Value = "  -4.60%  "

Set-TextValue $ws.Range("D33") "0.09265"
$ws.Range("E33").Value = "  -3.91%  "

Set-TextValue $ws.Range("D34") "1.458"
$ws.Range("E34").Value = "  +0.20%  "

Set-TextValue $ws.Range("D35") "3.564"
$ws.Range("E35").Value = "  -3.63%  "

Set-TextValue $ws.Range("D36") "5.288"
$ws.Range("E36").Value = "  -5.92%  "

Set-TextValue $ws.Range("D37") "0.02220"
$ws.Range("E37").Value = "  -4.91%  "

Set-TextValue $ws.Range("D38") "0.05972"
$ws.Range("E38").Value = "  -4.22%  "

Set-TextValue $ws.Range("D39") "1.205"
$ws.Range("E39").Value = "  -4.15%  "

Set-TextValue $ws.Range("D40") "8.278"
$ws.Range("E40").Value = "  -9.89%  "

Set-TextValue $ws.Range("D41") "0.9999"
$ws.Range("E41").Value = "  -1.17%  "

Set-TextValue $ws.Range("D42") "0.5893"
$ws.Range("E42").Value = "  -5.18%  "

Set-TextValue $ws.Range("D43") "0.1880"
$ws.Range("E43").Value = "  -1.78%  "

Set-TextValue $ws.Range("D44") "10.08"
$ws.Range("E44").Value = "  -9.79%  "

Set-TextValue $ws.Range("D45") "1.262"
$ws.Range("E45").Value = "  -6.65%  "

Set-TextValue $ws.Range("D46") "0.5625"
$ws.Range("E46").Value = "  -5.45%  "

Set-TextValue $ws.Range("D47") "11.94"
$ws.Range("E47").Value = "  -8.41%  "

Set-TextValue $ws.Range("D50") "0.06804"
$ws.Range("E50").Value = "  -0.22%  "

Set-TextValue $ws.Range("D51") "108.29"
$ws.Range("E51").Value = "  -2.70%  "

Write-Host "Applied cryptos update"
